$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.549.90'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.081.38'
$ws.Range('E3').Value = '  -2.38%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.92'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.81'
$ws.Range('E6').Value = '  +3.83%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.561'
$ws.Range('E8').Value = '  +5.70%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.074.28'
$ws.Range('E9').Value = '  -2.33%  '
$ws.Range('E10').Value = '  -2.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.87'
$ws.Range('E11').Value = '  -0.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.463'
$ws.Range('E12').Value = '  +0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000242'
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.54'
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('B15').Value = 'TRON'
$ws.Range('C15').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.119'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.592.28'
$ws.Range('E16').Value = '  -2.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.19'
$ws.Range('E17').Value = '  -0.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '63.507.77'
$ws.Range('E18').Value = '  -1.06%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.079.99'
$ws.Range('E19').Value = '  -2.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '479.49'
$ws.Range('E20').Value = '  +2.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.71'
$ws.Range('E21').Value = '  +1.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.721'
$ws.Range('E22').Value = '  -2.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.57'
$ws.Range('E23').Value = '  +0.61%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.18'
$ws.Range('E24').Value = '  +0.59%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.36'
$ws.Range('E25').Value = '  +0.63%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '81.94'
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('E27').Value = '  +0.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.74'
$ws.Range('E28').Value = '  +1.18%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.68'
$ws.Range('E29').Value = '  -1.92%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.29'
$ws.Range('E30').Value = '  -1.95%  '
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.20'
$ws.Range('E32').Value = '  -2.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.45'
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0₃0851'
$ws.Range('E35').Value = '  +0.62%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.06'
$ws.Range('E36').Value = '  -1.62%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.40'
$ws.Range('E37').Value = '  +3.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.14'
$ws.Range('E38').Value = '  -1.39%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.23'
$ws.Range('E39').Value = '  -4.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '9.37'
$ws.Range('E40').Value = '  +1.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '50.63'
$ws.Range('E41').Value = '  -2.19%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '443.56'
$ws.Range('E42').Value = '  -4.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.285'
$ws.Range('E43').Value = '  -4.70%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.815.88'
$ws.Range('E45').Value = '  -4.01%  '
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '39.28'
$ws.Range('E47').Value = '  -2.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.09'
$ws.Range('E48').Value = '  +0.39%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '25.24'
$ws.Range('E50').Value = '  +3.39%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.113'
$ws.Range('E51').Value = '  +1.78%  '
